$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add columns L and M (copy formatting from K1, which already has the right style) ---
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

# --- Cells whose new text is all-digits (e.g. "0130", "1000", "325") would be silently
#     reinterpreted as numbers (and lose leading zeros) by a plain .Value assignment,
#     exactly like in real Excel. Force them to Text first so they land as shared strings. ---
$textCells = @("G2", "K2", "G3", "K3", "G4", "K4", "J5", "J6")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- write the new / rearranged values ---
# Row 2
$ws.Range("G2").Value = "1000"
$ws.Range("H2").Value = "Beginning Arabic II MTWRF"
$ws.Range("I2").Value = "1050am"
$ws.Range("J2").Value = "KING"
$ws.Range("K2").Value = "325"
$ws.Range("L2").Value = "Al"
$ws.Range("M2").Value = "Raba'a Basem"

# Row 3
$ws.Range("G3").Value = "0130"
$ws.Range("H3").Value = "Language & Society Arab World TR"
$ws.Range("I3").Value = "0245pm"
$ws.Range("J3").Value = "KING"
$ws.Range("K3").Value = "325"
$ws.Range("L3").Value = "Al"
$ws.Range("M3").Value = "Raba'a Basem"

# Row 4
$ws.Range("G4").Value = "1100"
$ws.Range("H4").Value = "Intermediate Arabic II MTWRF"
$ws.Range("I4").Value = "1150am"
$ws.Range("J4").Value = "KING"
$ws.Range("K4").Value = "325"
$ws.Range("L4").Value = "Al"
$ws.Range("M4").Value = "Raba'a Basem"

# Row 5
$ws.Range("G5").Value = "II"
$ws.Range("H5").Value = "Intermediate Arabic"
$ws.Range("I5").Value = "MTWF"
$ws.Range("J5").Value = "0230"
$ws.Range("K5").Value = "0320pm"
$ws.Range("L5").Value = "TBA"
$ws.Range("M5").Value = "Hamilton Elizabeth"

# Row 6
$ws.Range("G6").Value = "Arabic"
$ws.Range("H6").Value = "Advanced"
$ws.Range("I6").Value = "MTWR"
$ws.Range("J6").Value = "1000"
$ws.Range("K6").Value = "1050am"
$ws.Range("L6").Value = "TBA"
$ws.Range("M6").Value = "Hamilton Elizabeth"

# Row 7
$ws.Range("G7").ClearContents() | Out-Null
$ws.Range("H7").Value = "Private Reading "
$ws.Range("I7").Value = "Full"
$ws.Range("J7").Value = "TBA"
$ws.Range("K7").Value = "TBA"
$ws.Range("L7").Value = "Al"
$ws.Range("M7").Value = "Raba'a Basem"

# Row 8
$ws.Range("G8").ClearContents() | Out-Null
$ws.Range("H8").Value = "Private Reading "
$ws.Range("I8").Value = "Half"
$ws.Range("J8").Value = "TBA"
$ws.Range("K8").Value = "TBA"
$ws.Range("L8").Value = "Al"
$ws.Range("M8").Value = "Raba'a Basem"

# --- now that the text is safely stored, drop the temporary Text number-format again so
#     the cells end up with no explicit style, just like the rest of the sheet ---
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats() | Out-Null
}
